$d = $word.ActiveDocument

# Title paragraph: "Answers: Laws of indices" (was split across 7 runs)
$d.Content.Find.Execute("Answers: Laws of indices", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Answers: Laws of indices", 2)

# Author paragraph: "Isabella Lewis, Akshat Srivastava" (was split across 7 runs)
$d.Content.Find.Execute("Isabella Lewis, Akshat Srivastava", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Isabella Lewis, Akshat Srivastava", 2)

# Abstract paragraph: "Answers to questions relating to using laws of indices." (was split across 15 runs)
$d.Content.Find.Execute("Answers to questions relating to using laws of indices.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Answers to questions relating to using laws of indices.", 2)
